$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row needs to be inserted at row 274 (pushing the
# existing rows 274-391 down to 275-392, growing the sheet from 391 to 392
# data rows). Insert a blank row first so everything below shifts down.
$ws.Rows(274).Insert()

# Populate the newly inserted row 274 with the new reading. The
# market/region/product/category boilerplate columns (A,B,C,E,F,G,H,I,J,R)
# match every other row for this product subset.
$ws.Range("A274").Value = 10
$ws.Range("B274").Value = "Vega Modelo de Temuco"
$ws.Range("C274").Value = "La Araucanía"
$ws.Range("D274").Value = 44992
$ws.Range("E274").Value = 9
$ws.Range("F274").Value = "Fruta"
$ws.Range("G274").Value = 100102
$ws.Range("H274").Value = "Cítricos"
$ws.Range("I274").Value = 100102006
$ws.Range("J274").Value = "Pomelo"
$ws.Range("K274").Value = "Start Ruby"
$ws.Range("L274").Value = "Primera"
$ws.Range("M274").Value = 80
$ws.Range("N274").Value = 14000
$ws.Range("O274").Value = 14000
$ws.Range("P274").Value = 14000
$ws.Range("Q274").Value = "$/bandeja 15 kilos granel"
$ws.Range("R274").Value = "Región de O'Higgins"
$ws.Range("S274").Value = 933
$ws.Range("T274").Value = 15
